$d = $word.ActiveDocument

# Locate the relevant paragraphs by content instead of hard-coded index,
# so the script is resilient to minor shifts elsewhere in the document.
$idxEsto = -1
$idxStar = -1
$n = $d.Paragraphs.Count
for ($i = 1; $i -le $n; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($t -like "*Esto pasa en tablas intermedias*") {
        $idxEsto = $i
    }
    if ($t.Trim() -eq "*") {
        $idxStar = $i
    }
}

# --- Edit 1 -----------------------------------------------------------
# The "Esto pasa en tablas intermedias..." paragraph currently holds a
# _GoBack bookmark BEFORE its run, and is followed by an empty paragraph.
# Rebuild both as:
#   - same paragraph, but with the run now appearing BEFORE the bookmark
#   - a brand-new paragraph (no paraId) with the 22/01/2020 note, ending
#     with the _GoBack bookmark
$pEsto = $d.Paragraphs.Item($idxEsto)
$pEmpty = $d.Paragraphs.Item($idxEsto + 1)
$combined = $d.Range($pEsto.Range.Start, $pEmpty.Range.End)

$xml1 = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">
<w:body>
<w:p w14:paraId="04788299" w14:textId="0D63B419" w:rsidR="00201068" w:rsidRDefault="00201068" w:rsidP="00AC2CE4">
<w:pPr><w:pStyle w:val="Prrafodelista"/><w:ind w:left="408"/></w:pPr>
<w:r><w:t>Esto pasa en tablas intermedias que deshacen la relación Many to Many. Por tanto en esas tablas solamente debemos guardar Id’s de los objetos que se relacionan y no el objeto entero.</w:t></w:r>
</w:p>
<w:p>
<w:r><w:t xml:space="preserve">22/01/2020 : Se crea proyecto web con ASPNetCore web api. Se genera nueva bbdd con SqlExpres (bbdd en : </w:t></w:r>
<w:r><w:t>C:\Program Files\Microsoft SQL Server\MSSQL14.SQLEXPRESS\MSSQL\DATA</w:t></w:r>
<w:r><w:t xml:space="preserve">), bbdd se llama </w:t></w:r>
<w:r><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/><w:color w:val="A31515"/><w:sz w:val="19"/><w:szCs w:val="19"/></w:rPr><w:t>WebAcademyDb</w:t></w:r>
<w:r><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/><w:color w:val="A31515"/><w:sz w:val="19"/><w:szCs w:val="19"/></w:rPr><w:t>.mdf. Sigue estando la bbdd de sqlite en el proyecto wpf.</w:t></w:r>
<w:bookmarkStart w:id="0" w:name="_GoBack"/>
<w:bookmarkEnd w:id="0"/>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@

$combined.InsertXML($xml1) | Out-Null

# --- Edit 2 -----------------------------------------------------------
# Paragraph holding only "*" loses its run, becoming an empty paragraph
# (keeps its own pPr/paraId).
$pStar = $d.Paragraphs.Item($idxStar)
$rStar = $pStar.Range

$xml2 = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">
<w:body>
<w:p w14:paraId="18779948" w14:textId="193CD556" w:rsidR="00201068" w:rsidRDefault="00201068" w:rsidP="00E01B32">
<w:pPr><w:pStyle w:val="Prrafodelista"/><w:ind w:left="408"/></w:pPr>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@

$rStar.InsertXML($xml2) | Out-Null
